# Scheduled-runner refresh of the FFXIV leve-profit calculations.
# Updates the currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) for the rows whose market-board snapshot changed, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 383.4
$ws.Range("I49").Value = 308.5
$ws.Range("J49").Value = 433.33334
$ws.Range("K49").Value = 925.5
$ws.Range("L49").Value = 1300.00002
$ws.Range("M49").Value = -789.5
$ws.Range("N49").Value = -1572.00002

$ws.Range("H69").Value = 5685339
$ws.Range("J69").Value = 3761.0527
$ws.Range("L69").Value = 11283.1581
$ws.Range("N69").Value = -13031.1581

$ws.Range("H72").Value = 5685339
$ws.Range("J72").Value = 3761.0527
$ws.Range("L72").Value = 33849.4743
$ws.Range("N72").Value = -42585.4743

$ws.Range("H98").Value = 530
$ws.Range("I98").Value = 530
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 530
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 968
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 530
$ws.Range("I122").Value = 530
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1590
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 860
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 1648.5294
$ws.Range("I137").Value = 1555.8148
$ws.Range("J137").Value = 1752.8334
$ws.Range("K137").Value = 4667.4444
$ws.Range("L137").Value = 5258.5002
$ws.Range("M137").Value = -2117.4444
$ws.Range("N137").Value = -10358.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8166.116
$ws.Range("I32").Value = 2709.9048
$ws.Range("J32").Value = 23111.39
$ws.Range("K32").Value = 2709.9048
$ws.Range("L32").Value = 23111.39
$ws.Range("M32").Value = -2422.9048
$ws.Range("N32").Value = -23685.39

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 110.111115
$ws.Range("I80").Value = 94.75
$ws.Range("J80").Value = 122.4
$ws.Range("K80").Value = 94.75
$ws.Range("L80").Value = 122.4
$ws.Range("M80").Value = 903.25
$ws.Range("N80").Value = -2118.4

$ws.Range("H83").Value = 110.111115
$ws.Range("I83").Value = 94.75
$ws.Range("J83").Value = 122.4
$ws.Range("K83").Value = 473.75
$ws.Range("L83").Value = 612
$ws.Range("M83").Value = 4518.25
$ws.Range("N83").Value = -10596

$ws.Range("H94").Value = 31435.857
$ws.Range("I94").Value = 2502.25
$ws.Range("J94").Value = 70014
$ws.Range("K94").Value = 2502.25
$ws.Range("L94").Value = 70014
$ws.Range("M94").Value = -2051.25
$ws.Range("N94").Value = -70916

$ws.Range("H107").Value = 3158.875
$ws.Range("I107").Value = 2881.8333
$ws.Range("J107").Value = 3990
$ws.Range("K107").Value = 2881.8333
$ws.Range("L107").Value = 3990
$ws.Range("M107").Value = -961.8332999999998
$ws.Range("N107").Value = -7830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372

$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864

$ws.Range("H132").Value = 2087.1875
$ws.Range("I132").Value = 1121.3478
$ws.Range("J132").Value = 4555.4443
$ws.Range("K132").Value = 3364.0434
$ws.Range("L132").Value = 13666.3329
$ws.Range("M132").Value = -834.0434
$ws.Range("N132").Value = -18726.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 111122570
$ws.Range("I50").Value = 500000100
$ws.Range("J50").Value = 14698.571
$ws.Range("K50").Value = 1500000300
$ws.Range("L50").Value = 44095.713
$ws.Range("M50").Value = -1499999819
$ws.Range("N50").Value = -45057.713

$ws.Range("H53").Value = 111122570
$ws.Range("I53").Value = 500000100
$ws.Range("J53").Value = 14698.571
$ws.Range("K53").Value = 1500000300
$ws.Range("L53").Value = 44095.713
$ws.Range("M53").Value = -1499999819
$ws.Range("N53").Value = -45057.713

$ws.Range("H131").Value = 1191353
$ws.Range("J131").Value = 969.0789
$ws.Range("L131").Value = 2907.2367
$ws.Range("N131").Value = -12987.2367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.375
$ws.Range("I2").Value = 46.2
$ws.Range("J2").Value = 246.66667
$ws.Range("K2").Value = 46.2
$ws.Range("L2").Value = 246.66667
$ws.Range("M2").Value = 66.8
$ws.Range("N2").Value = -472.66667

$ws.Range("H122").Value = 1773.2858
$ws.Range("I122").Value = 1868.9166
$ws.Range("K122").Value = 5606.7498
$ws.Range("M122").Value = -3156.7498

$ws.Range("H125").Value = 60800
$ws.Range("J125").Value = 60800
$ws.Range("L125").Value = 60800
$ws.Range("N125").Value = -65720

$ws.Range("H132").Value = 2209.077
$ws.Range("I132").Value = 1706.5264
$ws.Range("J132").Value = 3573.1428
$ws.Range("K132").Value = 5119.5792
$ws.Range("L132").Value = 10719.4284
$ws.Range("M132").Value = -2589.5792
$ws.Range("N132").Value = -15779.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4008.2856
$ws.Range("I7").Value = 3920.1
$ws.Range("J7").Value = 4228.75
$ws.Range("K7").Value = 3920.1
$ws.Range("L7").Value = 4228.75
$ws.Range("M7").Value = -3808.1
$ws.Range("N7").Value = -4452.75

$ws.Range("H16").Value = 1452.7273
$ws.Range("I16").Value = 1602.2222
$ws.Range("J16").Value = 780
$ws.Range("K16").Value = 1602.2222
$ws.Range("L16").Value = 780
$ws.Range("M16").Value = -1432.2222
$ws.Range("N16").Value = -1120

$ws.Range("H36").Value = 35000
$ws.Range("J36").Value = 35000
$ws.Range("L36").Value = 35000
$ws.Range("N36").Value = -36124

$ws.Range("H40").Value = 3914.3704
$ws.Range("I40").Value = 3761.8333
$ws.Range("J40").Value = 4219.4443
$ws.Range("K40").Value = 3761.8333
$ws.Range("L40").Value = 4219.4443
$ws.Range("M40").Value = -3625.8333
$ws.Range("N40").Value = -4491.4443

$ws.Range("H46").Value = 1399.1333
$ws.Range("I46").Value = 1212.4546
$ws.Range("J46").Value = 1912.5
$ws.Range("K46").Value = 1212.4546
$ws.Range("L46").Value = 1912.5
$ws.Range("M46").Value = -1024.4546
$ws.Range("N46").Value = -2288.5

$ws.Range("H55").Value = 246.10527
$ws.Range("I55").Value = 215
$ws.Range("J55").Value = 268.72726
$ws.Range("K55").Value = 215
$ws.Range("L55").Value = 268.72726
$ws.Range("M55").Value = -42
$ws.Range("N55").Value = -614.72726

$ws.Range("H126").Value = 4008.2856
$ws.Range("I126").Value = 3920.1
$ws.Range("J126").Value = 4228.75
$ws.Range("K126").Value = 11760.3
$ws.Range("L126").Value = 12686.25
$ws.Range("M126").Value = -9290.299999999999
$ws.Range("N126").Value = -17626.25

$ws.Range("H136").Value = 6585382.5
$ws.Range("J136").Value = 1282.2727
$ws.Range("L136").Value = 3846.8181
$ws.Range("N136").Value = -8946.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2933.4285
$ws.Range("I126").Value = 3255.6667
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 9767.000100000001
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -7297.000100000001
$ws.Range("N126").Value = -7940

$ws.Range("H136").Value = 7358.684
$ws.Range("I136").Value = 1443.75
$ws.Range("J136").Value = 11660.454
$ws.Range("K136").Value = 4331.25
$ws.Range("L136").Value = 34981.362
$ws.Range("M136").Value = -1781.25
$ws.Range("N136").Value = -40081.362
